$wb = $excel.ActiveWorkbook

# Sheet: Operating_cost_per_technology
$ws = $wb.Worksheets.Item("Operating_cost_per_technology")
$ws.Range("A1").Value = "Gas_CHP"
$ws.Range("B1").Value = 19410.459383650035
$ws.Range("A2").Value = "Gas_boiler"
$ws.Range("B2").Value = 113134.25260950485
$ws.Range("A3").Value = "Grid"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "Heat_pump"
$ws.Range("B4").Value = 0
$ws.Range("A5").Value = "Solar_PV"
$ws.Range("B5").Value = 0
$ws.Range("A6").Value = "Solar_thermal"
$ws.Range("B6").Value = 0
$ws.Range("A7:B7").Delete()

# Sheet: Maintenance_cost_per_technology
$ws = $wb.Worksheets.Item("Maintenance_cost_per_technology")
$ws.Range("A1").Value = "Gas_CHP"
$ws.Range("B1").Value = 3709.3387882160555
$ws.Range("A2").Value = "Gas_boiler"
$ws.Range("B2").Value = 11816.244161437475
$ws.Range("A3").Value = "Grid"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "Heat_pump"
$ws.Range("B4").Value = 0
$ws.Range("A5").Value = "Solar_PV"
$ws.Range("B5").Value = 14290.406901269387
$ws.Range("A6").Value = "Solar_thermal"
$ws.Range("B6").Value = 0
$ws.Range("A7:B7").Delete()

# Sheet: Capital_cost_per_technology
$ws = $wb.Worksheets.Item("Capital_cost_per_technology")
$ws.Range("A1").Value = "Gas_CHP"
$ws.Range("B1").Value = 3379.3738649993311
$ws.Range("A2").Value = "Gas_boiler"
$ws.Range("B2").Value = 36533.170275829085
$ws.Range("A3").Value = "Grid"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "Heat_pump"
$ws.Range("B4").Value = 0
$ws.Range("A5").Value = "Solar_PV"
$ws.Range("B5").Value = 34710.374943464696
$ws.Range("A6").Value = "Solar_thermal"
$ws.Range("B6").Value = 0
$ws.Range("A7:B7").Delete()

# Sheet: Total_cost_per_technology
$ws = $wb.Worksheets.Item("Total_cost_per_technology")
$ws.Range("A1").Value = "Gas_CHP"
$ws.Range("B1").Value = 26499.172036865421
$ws.Range("A2").Value = "Gas_boiler"
$ws.Range("B2").Value = 161483.66704677141
$ws.Range("A3").Value = "Grid"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "Heat_pump"
$ws.Range("B4").Value = 0
$ws.Range("A5").Value = "Solar_PV"
$ws.Range("B5").Value = 49000.781844734083
$ws.Range("A6").Value = "Solar_thermal"
$ws.Range("B6").Value = 0
$ws.Range("A7:B7").Delete()

# Sheet: Operating_cost_grid
$ws = $wb.Worksheets.Item("Operating_cost_grid")
$ws.Range("A1").Value = 220939.27557284833

# Sheet: Total_cost_grid
$ws = $wb.Worksheets.Item("Total_cost_grid")
$ws.Range("A1").Value = 220939.27557284833

# Sheet: Capital_cost_per_storage
$ws = $wb.Worksheets.Item("Capital_cost_per_storage")
$ws.Range("A1").Value = "Elec"
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = "Heat"
$ws.Range("B2").Value = 3656.913220251934

# Sheet: Total_cost_per_storage
$ws = $wb.Worksheets.Item("Total_cost_per_storage")
$ws.Range("A1").Value = "Elec"
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = "Heat"
$ws.Range("B2").Value = 3656.913220251934

# Sheet: Income_via_exports
$ws = $wb.Worksheets.Item("Income_via_exports")
$ws.Range("A1").Value = 4417.7851897133714
